# Auto-generated Excel COM-interop edit script
# Applies the 28.01.2021 daily update to the Slovakia Covid DailyStats sheet:
#  - updates AgTests (H) / AgPosit (I) retro-corrections for rows 271-328
#  - appends a new row 329 for date serial 44223 (2021-01-27)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to existing AgTests (H) / AgPosit (I) values ---
$ws.Range("H271").Value = 42644
$ws.Range("I271").Value = 1624
$ws.Range("H272").Value = 31038
$ws.Range("H273").Value = 27170
$ws.Range("I273").Value = 1367
$ws.Range("H274").Value = 28505
$ws.Range("I274").Value = 1349
$ws.Range("H275").Value = 28935
$ws.Range("I275").Value = 1243
$ws.Range("H278").Value = 30137
$ws.Range("I278").Value = 2108
$ws.Range("H279").Value = 43499
$ws.Range("H280").Value = 35723
$ws.Range("I280").Value = 2409
$ws.Range("H281").Value = 45351
$ws.Range("I281").Value = 3257
$ws.Range("H282").Value = 46962
$ws.Range("I282").Value = 2862
$ws.Range("H285").Value = 41051
$ws.Range("I285").Value = 3446
$ws.Range("H286").Value = 54367
$ws.Range("I286").Value = 4261
$ws.Range("H287").Value = 57813
$ws.Range("I287").Value = 3927
$ws.Range("H288").Value = 56278
$ws.Range("I288").Value = 3961
$ws.Range("H289").Value = 64611
$ws.Range("I289").Value = 3707
$ws.Range("H292").Value = 81543
$ws.Range("I292").Value = 7210
$ws.Range("H293").Value = 81994
$ws.Range("I293").Value = 5800
$ws.Range("H294").Value = 91081
$ws.Range("I294").Value = 5044
$ws.Range("H299").Value = 64312
$ws.Range("I299").Value = 6731
$ws.Range("H300").Value = 70686
$ws.Range("I300").Value = 6929
$ws.Range("H301").Value = 69831
$ws.Range("I301").Value = 5566
$ws.Range("H302").Value = 72786
$ws.Range("H306").Value = 70962
$ws.Range("I306").Value = 7186
$ws.Range("H307").Value = 73503
$ws.Range("I307").Value = 6330
$ws.Range("H309").Value = 57439
$ws.Range("I309").Value = 3971
$ws.Range("H310").Value = 91188
$ws.Range("I310").Value = 5198
$ws.Range("H313").Value = 73159
$ws.Range("I313").Value = 3560
$ws.Range("H314").Value = 65283
$ws.Range("I314").Value = 3362
$ws.Range("H315").Value = 66489
$ws.Range("I315").Value = 3018
$ws.Range("H316").Value = 49147
$ws.Range("I316").Value = 2281
$ws.Range("H317").Value = 61819
$ws.Range("I317").Value = 2153
$ws.Range("H320").Value = 86520
$ws.Range("I320").Value = 3912
$ws.Range("H321").Value = 90442
$ws.Range("I321").Value = 2799
$ws.Range("H322").Value = 104742
$ws.Range("I322").Value = 2286
$ws.Range("H323").Value = 149944
$ws.Range("I323").Value = 2328
$ws.Range("H324").Value = 231451
$ws.Range("I324").Value = 2669
$ws.Range("H325").Value = 667595
$ws.Range("I325").Value = 5477
$ws.Range("H326").Value = 402981
$ws.Range("I326").Value = 3517
$ws.Range("H327").Value = 253197
$ws.Range("I327").Value = 3772
$ws.Range("H328").Value = 188092
$ws.Range("I328").Value = 2734

# --- New row 329: daily update for 2021-01-27 (serial 44223) ---
$ws.Range("A329").Value = 44223
$ws.Range("B329").Value = 243427
$ws.Range("C329").Value = 205247
$ws.Range("D329").Value = 33769
$ws.Range("E329").Value = 9811
$ws.Range("F329").Value = 2035
$ws.Range("G329").Value = 4411
$ws.Range("H329").Value = 78633
$ws.Range("I329").Value = 1779
